$d = $word.ActiveDocument

# Locate the exact text "Oklahoma City," (the SenderInfo line that reads
# "Oklahoma City, OK") without mutating anything yet.
$found = $d.Content
$found.Find.ClearFormatting()
$found.Find.Execute("Oklahoma City,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found.Find.Found) {
    throw "Could not find target text 'Oklahoma City,'"
}

$matchStart = $found.Start
$matchEnd = $found.End

# Find the paragraph that contains the match so we can rewrite its whole
# (non-paragraph-mark) content in one shot -- this engine folds any runs
# left dangling next to a partial edit back together, so the safest way to
# hit the exact run layout from the diff is to replace the complete
# paragraph body rather than only the matched substring.
$paraCount = $d.Paragraphs.Count
$pStart = -1
$pEnd = -1
$pIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $matchStart -and $p.Range.End -ge $matchEnd) {
        $pStart = $p.Range.Start
        $pEnd = $p.Range.End
        $pIndex = $i
        break
    }
}
if ($pStart -eq -1) {
    throw "Could not locate enclosing paragraph for match"
}

# Paragraph.Range.End includes the trailing paragraph-mark character, so
# drop the last position to get the pure text span.
$bodyStart = $pStart
$bodyEnd = $pEnd - 1

$body = $d.Range($bodyStart, $bodyEnd)
Write-Output "Original paragraph text: [$($body.Text)]"
$body.Delete()

# Re-insert the full paragraph content: the expanded city-name run, a new
# run holding just the trailing comma (both Verdana/sz-20, matching the
# diff), then the untouched " OK" tail reconstructed as its original two
# runs.
$insertion = $d.Range($bodyStart, $bodyStart)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00821BC8"><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr><w:t>Oklahoma City Metropolitan Area</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr><w:t>,</w:t></w:r><w:r w:rsidR="009E0C87"><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00821BC8"><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana"/><w:sz w:val="20"/></w:rPr><w:t>OK</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertion.InsertXML($xml)

$after = $d.Paragraphs.Item($pIndex).Range.Text
Write-Output "Updated paragraph text: [$after]"
